$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160, shifting existing rows 160:228 down to 161:229.
$ws.Rows(160).Insert()

# Populate the newly inserted row 160 with the new weekly price-record data.
$ws.Range("A160").Value = 8
$ws.Range("B160").Value = "Terminal La Palmera de La Serena"
$ws.Range("C160").Value = "Coquimbo"
$ws.Range("D160").Value = 44704
$ws.Range("D160").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E160").Value = 4
$ws.Range("F160").Value = 100112031
$ws.Range("G160").Value = "Poroto verde"
$ws.Range("H160").Value = "Magnum"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 400
$ws.Range("K160").Value = 26000
$ws.Range("L160").Value = 27000
$ws.Range("M160").Value = 26500
$ws.Range("N160").Value = "$/malla 25 kilos"
$ws.Range("O160").Value = "Perú"
$ws.Range("P160").Value = 1060
$ws.Range("Q160").Value = 25
$ws.Range("R160").Value = "Hortaliza"
